# Weekly update: insert a new daily price record as row 81 in the
# "Hortaliza, Vega Modelo de Temuco - Berenjena" sheet. All the rows that
# used to be 81..188 shift down to 82..189 (handled automatically by the
# row Insert below), and the new row 81 is populated with the latest
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81, pushing existing rows 81-188 down
# to 82-189.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new data point.
$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "Vega Modelo de Temuco"
$ws.Range("C81").Value = "La Araucanía"
$ws.Range("D81").Value = 44483
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = 100112001
$ws.Range("G81").Value = "Berenjena"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 150
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = 10667
$ws.Range("N81").Value = "`$/caja 60 unidades"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 178
$ws.Range("Q81").Value = 60
$ws.Range("R81").Value = "Hortaliza"
